$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.310.68"
$ws.Range("E2").Value = "  -0.29%  "
$ws.Range("D3").Value = "1.788.86"
$ws.Range("E3").Value = "  -1.57%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").Value = "'334.92"
$ws.Range("E5").Value = "  -2.86%  "
$ws.Range("D6").Value = "'0.9994"
$ws.Range("E6").Value = "  -0.11%  "
$ws.Range("E7").Value = "  -1.68%  "
$ws.Range("D8").Value = "'0.3449"
$ws.Range("E8").Value = "  -2.45%  "
$ws.Range("D9").Value = "'48.42"
$ws.Range("E9").Value = "  -3.90%  "
$ws.Range("D10").Value = "'1.205"
$ws.Range("E10").Value = "  -3.34%  "
$ws.Range("D11").Value = "'0.07518"
$ws.Range("E11").Value = "  -3.60%  "
$ws.Range("D12").Value = "'0.9994"
$ws.Range("E12").Value = "  -0.21%  "
$ws.Range("D13").Value = "'21.96"
$ws.Range("E13").Value = "  -3.97%  "
$ws.Range("D14").Value = "'6.503"
$ws.Range("E14").Value = "  -2.72%  "
$ws.Range("D15").Value = "1.789.44"
$ws.Range("E15").Value = "  -1.44%  "
$ws.Range("D16").Value = "'7.116"
$ws.Range("E16").Value = "  -2.11%  "
$ws.Range("D17").Value = "'0.00001103"
$ws.Range("E17").Value = "  -2.91%  "
$ws.Range("D18").Value = "'0.06660"
$ws.Range("E18").Value = "  -1.71%  "
$ws.Range("D19").Value = "'84.12"
$ws.Range("E19").Value = "  -3.24%  "
$ws.Range("D20").Value = "'1.000"
$ws.Range("E20").Value = "  +0.01%  "
$ws.Range("D21").Value = "'6.659"
$ws.Range("E21").Value = "  +1.02%  "
$ws.Range("D22").Value = "'17.42"
$ws.Range("E22").Value = "  -3.27%  "
$ws.Range("D23").Value = "27.303.76"
$ws.Range("E23").Value = "  -0.27%  "
$ws.Range("D24").Value = "'12.44"
$ws.Range("E24").Value = "  -5.85%  "
$ws.Range("D25").Value = "'2.426"
$ws.Range("E25").Value = "  -1.75%  "
$ws.Range("D26").Value = "'1.505"
$ws.Range("E26").Value = "  -0.80%  "
$ws.Range("D27").Value = "'2.567"
$ws.Range("E27").Value = "  -6.95%  "
$ws.Range("D28").Value = "'21.45"
$ws.Range("E28").Value = "  -2.83%  "
$ws.Range("D29").Value = "'153.94"
$ws.Range("E29").Value = "  -0.26%  "
$ws.Range("D30").Value = "1.991.98"
$ws.Range("E30").Value = "  -1.27%  "
$ws.Range("D31").Value = "'134.44"
$ws.Range("E31").Value = "  -2.29%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "'6.150"
$ws.Range("E32").Value = "  -4.70%  "
$ws.Range("B33").Value = "HuobiToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D33").Value = "'4.024"
$ws.Range("E33").Value = "  -2.25%  "
$ws.Range("D34").Value = "'0.08707"
$ws.Range("E34").Value = "  -1.57%  "
$ws.Range("D35").Value = "'13.34"
$ws.Range("E35").Value = "  -4.49%  "
$ws.Range("D36").Value = "'1.663"
$ws.Range("E36").Value = "  -3.67%  "
$ws.Range("D37").Value = "'0.6991"
$ws.Range("E37").Value = "  -3.17%  "
$ws.Range("D38").Value = "'5.492"
$ws.Range("E38").Value = "  -3.33%  "
$ws.Range("D39").Value = "'0.2214"
$ws.Range("E39").Value = "  -3.25%  "
$ws.Range("D40").Value = "'8.873"
$ws.Range("E40").Value = "  -2.66%  "
$ws.Range("D41").Value = "'0.06357"
$ws.Range("E41").Value = "  -3.86%  "
$ws.Range("D42").Value = "'0.02348"
$ws.Range("E42").Value = "  -3.51%  "
$ws.Range("D43").Value = "'1.247"
$ws.Range("E43").Value = "  -1.66%  "
$ws.Range("D44").Value = "'14.43"
$ws.Range("E44").Value = "  -3.75%  "
$ws.Range("D45").Value = "'0.6568"
$ws.Range("E45").Value = "  -2.61%  "
$ws.Range("D46").Value = "'0.9986"
$ws.Range("E46").Value = "  -0.15%  "
$ws.Range("D47").Value = "'3.850"
$ws.Range("E47").Value = "  -3.23%  "
$ws.Range("D48").Value = "'2.158"
$ws.Range("E48").Value = "  -2.18%  "
$ws.Range("D49").Value = "'129.52"
$ws.Range("E49").Value = "  -3.09%  "
$ws.Range("D50").Value = "'0.07154"
$ws.Range("E50").Value = "  -2.96%  "
$ws.Range("D51").Value = "'79.51"
$ws.Range("E51").Value = "  -2.28%  "
